$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the analysis labels (Burnout_* / instabilidade_* -> shorter fem/masc style labels)
$labels = @{
    "Burnout_cinismo"      = "Cinismo"
    "Burnout_exaustao"     = "Exaustão"
    "Burnout_eficacia"     = "Eficacia"
    "instabilidade_equipe" = "Inst.Equipe"
    "instabilidade_tecnica"= "Inst.Tecnica"
    "instabilidade_tarefa" = "Inst.Tarefa"
}

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value2
    if ($null -ne $val -and $labels.ContainsKey([string]$val)) {
        $cell.Value = $labels[[string]$val]
    }
}
